$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E7"  = 15.33060000000001
    "C8"  = -12.30889999999999
    "C10" = -13.24
    "C12" = -10.6487
    "E15" = 16.1556
    "C18" = -11.62579999999999
    "E18" = 18.13400000000001
    "E20" = 15.97759999999999
    "E29" = 17.03750000000002
    "E30" = 15.5702
    "E31" = 15.94890000000001
    "C37" = -14.23269999999999
    "E40" = 17.07560000000001
    "E50" = 16.16239999999999
    "C55" = -13.24309999999999
    "C68" = -11.1096
    "E68" = 17.00080000000001
    "E76" = 16.10349999999999
    "C77" = -11.8159
    "C78" = -12.14240000000001
    "C81" = -12.9249
    "C82" = -11.95559999999999
    "E87" = 16.221
    "E88" = 16.33530000000001
    "E96" = 16.2061
    "E98" = 15.5965
    "E101" = 16.75540000000001
    "E102" = 16.6086
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
